# 6.2.1 indicator sheet: add 2014 and 2023 data columns alongside the
# existing 2018 column, and update the footnote to mention all three
# MICS survey years.
#
# Shape of the change:
#   - Column D currently holds the single "2018" series (header + 23 data
#     rows, rows 4-30; row 3 is just a thin border row and row 31 is the
#     footnote row, neither of which gets new columns).
#   - We copy that whole D4:D30 block (values + formatting) into a new
#     column E, so the 2018 numbers end up preserved one column to the
#     right.
#   - Column D is then overwritten with the new 2014 figures, and a new
#     column F is filled with the 2023 figures using D's formatting.
#   - The footnote row (31) text is updated to reference all three survey
#     years instead of just 2018.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the existing "2018" column (D4:D30) into the new column E ---
$ws.Range("D4:D30").Copy()
$ws.Range("E4:E30").PasteSpecial(-4104)   # xlPasteAll (values + formats)
$excel.CutCopyMode = $false

# --- 2. Give the new column F the same formatting as D (values set below) ---
$ws.Range("D4:D30").Copy()
$ws.Range("F4:F30").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Header row: years 2014 / 2018 / 2023 ---
$ws.Range("D4").Value = 2014
$ws.Range("F4").Value = 2023

# --- 4. Data rows: new 2014 (D) and 2023 (F) figures ---
$ws.Range("D5").Value = 94.5
$ws.Range("F5").Value = 97.6

$ws.Range("D7").Value = 96.8
$ws.Range("F7").Value = 96.7

$ws.Range("D8").Value = 93
$ws.Range("F8").Value = 98

$ws.Range("D10").Value = 91.8
$ws.Range("F10").Value = 97

$ws.Range("D11").Value = 93
$ws.Range("F11").Value = 98.1

$ws.Range("D12").Value = 95.2
$ws.Range("F12").Value = 98.7

$ws.Range("D13").Value = 97.3
$ws.Range("F13").Value = 99.5

$ws.Range("D14").Value = 94.3
$ws.Range("F14").Value = 97.3

$ws.Range("D15").Value = 91
$ws.Range("F15").Value = 97.1

$ws.Range("D16").Value = 92.7
$ws.Range("F16").Value = 99.1

$ws.Range("D17").Value = 99.1
$ws.Range("F17").Value = 97.8

$ws.Range("D18").Value = 92.4
$ws.Range("F18").Value = 90.4

$ws.Range("D20").Value = 95.1
$ws.Range("F20").Value = 91.1

$ws.Range("D21").Value = 90.6
$ws.Range("F21").Value = 97.5

$ws.Range("D22").Value = 93.6
$ws.Range("F22").Value = 97.6

$ws.Range("D23").Value = 95.6
$ws.Range("F23").Value = 97.9

$ws.Range("D24").Value = 97.4
$ws.Range("F24").Value = 98

$ws.Range("D26").Value = 88.2
$ws.Range("F26").Value = 95.1

$ws.Range("D27").Value = 93.4
$ws.Range("F27").Value = 97.7

$ws.Range("D28").Value = 95.1
$ws.Range("F28").Value = 98.8

$ws.Range("D29").Value = 95.8
$ws.Range("F29").Value = 99

$ws.Range("D30").Value = 97.9
$ws.Range("F30").Value = 97.6

# --- 5. Footnote row (31): mention 2014 / 2018 / 2023 instead of just 2018 ---
$ws.Range("A31").Value = " Көп көрсөткүчтүү кластердик изилдөөнүн маалыматтары боюнча, 2014- ж., 2018-ж., 2023-ж."
$ws.Range("B31").Value = "По данным кластерного обследования по многим показателям, 2014г., 2018г., 2023г."
$ws.Range("C31").Value = "According to Multiple Indicator Cluster Survey, 2014, 2018, 2023."

# --- 6. Row 1 height shrinks slightly now that the header wraps differently ---
$ws.Rows("1").RowHeight = 67.5

# --- 7. Put the selection back on A1 ---
$ws.Range("A1").Select()
